$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("B5").Value = 91829
$ws.Range("B6").Value = 79244
$ws.Range("B7").Value = 80349
$ws.Range("B8").Value = 79715
$ws.Range("B10").Value = 79245
$ws.Range("B11").Value = 79245
$ws.Range("B12").Value = 80349
$ws.Range("B13").Value = 91829
$ws.Range("B14").Value = 91829
$ws.Range("B15").Value = 91829
